$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.144.48'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.834.31'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''232.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4660'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('D8').Value = '''0.2709'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.10%  '
$ws.Range('E9').Value = '  -3.46%  '
$ws.Range('D10').Value = '1.833.06'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').Value = '''0.07399'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '''15.98'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = '''4.913'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('E14').Value = '  -4.25%  '
$ws.Range('D15').Value = '''0.6153'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.52%  '
$ws.Range('D16').Value = '30.063.74'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '''228.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = '''0.000007263'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.43%  '
$ws.Range('D20').Value = '''1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '''12.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.92%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.076.31'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''4.843'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.82%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '''5.808'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.65%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '''9.182'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''165.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.99%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''17.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.55%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''1.862'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '''0.1027'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''1.371'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '''4.066'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.34%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''3.777'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.14%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.04776'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''1.130'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.02%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.7069'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.715'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.01864'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.644'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''0.8892'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''1.926'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.22%  '
$ws.Range('D41').Value = '''104.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.87%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.477'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.3994'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.65%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''6.928'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''59.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.1185'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''8.548'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.05510'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''32.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.352'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.72%  '
